# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly scraped output, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Values for worksheet "展览" (rows -> new F value)
$sheet1Updates = @{
    "F3"  = 5180
    "F5"  = 7483
    "F7"  = 71
    "F9"  = 602
    "F11" = 31
    "F12" = 4331
    "F13" = 1771
    "F14" = 107
    "F15" = 109
    "F16" = 2930
    "F19" = 211
    "F20" = 506
    "F21" = 444
    "F22" = 463
    "F23" = 312
    "F24" = 104
    "F25" = 1697
    "F26" = 1194
    "F27" = 94
    "F28" = 1388
    "F29" = 109
    "F30" = 583
    "F35" = 106
    "F36" = 68
    "F37" = 2929
    "F38" = 708
    "F39" = 27
    "F40" = 90
    "F41" = 43
    "F42" = 45
}

# Values for worksheet "全部类型" (rows -> new F value)
$sheet4Updates = @{
    "F3"  = 5180
    "F5"  = 7483
    "F7"  = 71
    "F9"  = 602
    "F11" = 31
    "F12" = 4331
    "F13" = 1771
    "F14" = 107
    "F15" = 109
    "F16" = 2930
    "F19" = 211
    "F20" = 506
    "F21" = 444
    "F22" = 463
    "F24" = 312
    "F25" = 104
    "F26" = 1697
    "F27" = 1194
    "F28" = 94
    "F29" = 1388
    "F30" = 109
    "F31" = 583
    "F36" = 106
    "F37" = 68
    "F38" = 2929
    "F40" = 708
    "F41" = 27
    "F42" = 90
    "F43" = 43
    "F44" = 45
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
